$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.713.02'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.426.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.67%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '613.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.417.28'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.93%  '
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.194'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.95'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.559'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.25%  '
$ws.Range("E14").Value = '  -1.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.979.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.67%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.431.76'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.87%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.890.76'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '573.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.95%  '
$ws.Range("E20").Value = '  +0.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.840'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '95.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.58%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '32.46'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.77'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("B32").Value = 'Mantle'
$ws.Range("C32").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.27%  '
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '596.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.27%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0949'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0464'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '55.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.25%  '
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -13.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.220.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₃0671'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.39%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.76'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.25%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.292'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.36'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.14%  '
$ws.Range("E49").Value = '  -2.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.89%  '
